# Generate Report for Handoff
# - Refresh the "Latest HO Xliff Generate Date" / handoff-handback timestamps
#   for the files that are "Ready for handoff" (rows 7, 8, 11-14).
# - Stamp the "Priority" column ("ht") for those same rows on the zh-cn and
#   de-de sheets, which was previously blank.

$wb = $excel.ActiveWorkbook

$reportRows = 7,8,11,12,13,14

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

foreach ($r in $reportRows) {
    # Overview sheet: Latest HO Xliff Generate Date (column G)
    $wsOverview.Range("G$r").Value = "2016-09-04 12:25:13"

    # zh-cn sheet: Priority (E) + Latest Handback DateTime (H)
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-04 12:25:06"

    # de-de sheet: Priority (E) + Latest Handback DateTime (H)
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-04 12:25:13"
}
